$d = $word.ActiveDocument

# The document splits each "<id>...</id>" marker across three runs:
#   run1 "<id>"        (Courier New, color 7f6000, sz 18)
#   run2 "<label>"     (Arial, color 000000, sz 22)
#   run3 "</id>"       (Courier New, color 7f6000, sz 18)
# The edit merges these into a single run "<id><label></id>" that keeps
# run1's formatting (matching what the first run already carries).

function Merge-IdTag([string]$label) {
    $full = "<id>" + $label + "</id>"

    $findRange = $d.Content
    $ok = $findRange.Find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        return
    }

    $start = $findRange.Start
    $end = $findRange.End

    # "<id>" is always 4 characters; keep that first run untouched (it
    # already has the desired formatting) and fold everything after it
    # (the label + "</id>") back into that same run.
    $afterOpenTag = $start + 4

    $rest = $d.Range($afterOpenTag, $end)
    $rest.Delete()

    $insertionPoint = $d.Range($afterOpenTag, $afterOpenTag)
    $insertionPoint.InsertAfter($label + "</id>")
}

Merge-IdTag "p039r_1"
Merge-IdTag "p039r_2"
